# Applies the crypto price/volume/coin-name updates described in the commit.
# Price (column D) values that look numeric must be written with the cell
# pre-formatted as Text ("@") so Excel stores them as the original literal
# string (e.g. "586.53") instead of silently converting to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "61.889.27"
$ws.Range("E2").Value = "  -1.94%  "

# Row 3
$ws.Range("D3").Value = "2.966.54"
$ws.Range("E3").Value = "  -2.91%  "

# Row 4
$ws.Range("E4").Value = "  +0.21%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.53"
$ws.Range("E5").Value = "  -0.06%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.70"
$ws.Range("E6").Value = "  -6.23%  "

# Row 7
$ws.Range("E7").Value = "  +0.17%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.521"
$ws.Range("E8").Value = "  -3.02%  "

# Row 9
$ws.Range("D9").Value = "2.959.10"
$ws.Range("E9").Value = "  -3.15%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.145"
$ws.Range("E10").Value = "  -5.96%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.72"
$ws.Range("E11").Value = "  -2.38%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.460"
$ws.Range("E12").Value = "  +2.47%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000225"
$ws.Range("E13").Value = "  -4.12%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.05"
$ws.Range("E14").Value = "  -6.13%  "

# Row 15
$ws.Range("E15").Value = "  +1.62%  "

# Row 16
$ws.Range("D16").Value = "3.466.04"
$ws.Range("E16").Value = "  -2.64%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.98"
$ws.Range("E17").Value = "  -2.17%  "

# Row 18
$ws.Range("D18").Value = "61.896.45"
$ws.Range("E18").Value = "  -1.97%  "

# Row 19
$ws.Range("D19").Value = "2.987.66"
$ws.Range("E19").Value = "  -2.18%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "450.65"
$ws.Range("E20").Value = "  -5.79%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.79"
$ws.Range("E21").Value = "  -3.44%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.678"
$ws.Range("E22").Value = "  -3.96%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.30"
$ws.Range("E23").Value = "  -2.86%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.22"
$ws.Range("E24").Value = "  -1.14%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.10"
$ws.Range("E25").Value = "  -4.82%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.15"
$ws.Range("E26").Value = "  -10.65%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  -0.27%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.57"
$ws.Range("E28").Value = "  -9.05%  "

# Row 29
$ws.Range("E29").Value = "  +0.16%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.63"
$ws.Range("E30").Value = "  -1.59%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.84"
$ws.Range("E31").Value = "  -7.21%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.06"
$ws.Range("E32").Value = "  -6.23%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.31"
$ws.Range("E33").Value = "  -1.20%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.107"
$ws.Range("E34").Value = "  -3.55%  "

# Row 35
$ws.Range("B35").Value = "Mantle"
$ws.Range("C35").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.01"
$ws.Range("E35").Value = "  -4.37%  "

# Row 36
$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").Value = "0.0₃0783"
$ws.Range("E36").Value = "  -4.59%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.66"
$ws.Range("E37").Value = "  -4.29%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.07"
$ws.Range("E38").Value = "  -6.48%  "

# Row 39
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.03"
$ws.Range("E39").Value = "  -0.79%  "

# Row 40
$ws.Range("B40").Value = "Cosmos"
$ws.Range("C40").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.09"
$ws.Range("E40").Value = "  -1.64%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.120"
$ws.Range("E41").Value = "  +4.14%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.82"
$ws.Range("E42").Value = "  -13.86%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "389.43"
$ws.Range("E43").Value = "  -10.48%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0352"
$ws.Range("E44").Value = "  -2.52%  "

# Row 45
$ws.Range("D45").Value = "2.719.50"
$ws.Range("E45").Value = "  -3.85%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.263"
$ws.Range("E46").Value = "  -9.39%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "36.68"
$ws.Range("E47").Value = "  -4.06%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.25"
$ws.Range("E48").Value = "  +0.58%  "

# Row 49
$ws.Range("E49").Value = "  +0.09%  "

# Row 50
$ws.Range("E50").Value = "  -1.72%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.15"
$ws.Range("E51").Value = "  -2.16%  "
